# ERP-347 - As a Case Work Clerk I want the venue details to be populated
# on any correspondence documents that are created.
#
# Populates the Manchester / Glasgow Employment Tribunal venue details
# (address, telephone, fax, DX, email) that were previously placeholder
# values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Manchester tribunal ---------------------------------------------
$ws.Range("B3").Value = "Manchester Employment Tribunal, Alexandra House, 14-22 The Parsonage, Manchester, M3 2JA"
$ws.Range("B4").Value = "0161 833 6100"
$ws.Range("B5").Value = "0870 739 4433"
$ws.Range("B6").Value = "DX 743570"
$ws.Range("B7").Value = "Manchesteret@justice.gov.uk"

# --- Glasgow tribunal ---------------------------------------------
$ws.Range("B8").Value = "Eagle Building, 215 Bothwell Street, Glasgow, G2 7TS"
$ws.Range("B9").Value = "0141 204 0730"
$ws.Range("B10").Value = "01264 785 177"
$ws.Range("B11").Value = "DX 7435701"
$ws.Range("B12").Value = "glasgowet@justice.gov.uk"

# Update the display text of the existing mailto hyperlinks on B7/B12 in
# place (iterate the collection instead of indexing via Item(), which
# would append a brand-new, id-less hyperlink rather than editing the
# existing one).
foreach ($h in $ws.Hyperlinks) {
    if ($h.Address -eq "mailto:manchester@gmail.com") {
        $h.TextToDisplay = "Manchesteret@justice.gov.uk"
    } elseif ($h.Address -eq "mailto:glasgow@gmail.com") {
        $h.TextToDisplay = "glasgowet@justice.gov.uk"
    }
}
